$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update numeric values in row 2
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 0.07111099999999999
$ws.Range("N2").Value = 0.142222
$ws.Range("O2").Value = 0.02711460746047303
$ws.Range("P2").Value = 0.02678527708115022
$ws.Range("Q2").Value = 0.008143394683333333
$ws.Range("R2").Value = 0.0488603681
$ws.Range("S2").Value = 0.02711460746047303
$ws.Range("T2").Value = 0.02678527708115022

# Update numeric values in row 3
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("O3").Value = 0.0245904030281302
$ws.Range("P3").Value = 0.03643759694506741
$ws.Range("S3").Value = 0.0245904030281302
$ws.Range("T3").Value = 0.03643759694506741

# Update numeric values in row 4
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("M4").Value = 2.4870065
$ws.Range("N4").Value = 4.974013
$ws.Range("O4").Value = 0.9482949895113968
$ws.Range("P4").Value = 0.9367771259737823
$ws.Range("Q4").Value = 0.2848036943583334
$ws.Range("R4").Value = 1.70882216615
$ws.Range("S4").Value = 0.9482949895113968
$ws.Range("T4").Value = 0.9367771259737823

# Remove rows 5,6,7 (data for the second "Sending cluster" group was removed)
$ws.Range("A5:T7").EntireRow.Delete()
